$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 109

$ws.Range("H9").Value = 424.68182
$ws.Range("I9").Value = 202.38889
$ws.Range("J9").Value = 1425
$ws.Range("K9").Value = 202.38889
$ws.Range("L9").Value = 1425
$ws.Range("M9").Value = -33.38889
$ws.Range("N9").Value = -1763

$ws.Range("H12").Value = 237
$ws.Range("I12").Value = 237
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 237
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -67

$ws.Range("H19").Value = 391.15384
$ws.Range("I19").Value = 384.25
$ws.Range("J19").Value = 394.22223
$ws.Range("K19").Value = 384.25
$ws.Range("L19").Value = 394.22223
$ws.Range("M19").Value = -209.25
$ws.Range("N19").Value = -744.2222300000001

$ws.Range("H62").Value = 2051.3333
$ws.Range("I62").Value = 1961.6
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 1961.6
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = -1337.6
$ws.Range("N62").Value = -3748

$ws.Range("H65").Value = 2051.3333
$ws.Range("I65").Value = 1961.6
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 9808
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = -6688
$ws.Range("N65").Value = -18740

$ws.Range("H100").Value = 1337.762
$ws.Range("I100").Value = 1159
$ws.Range("J100").Value = 1909.8
$ws.Range("K100").Value = 1159
$ws.Range("L100").Value = 1909.8
$ws.Range("M100").Value = -618
$ws.Range("N100").Value = -2991.8

$ws.Range("H132").Value = 4194.7744
$ws.Range("I132").Value = 4298.956
$ws.Range("J132").Value = 1833.3334
$ws.Range("K132").Value = 12896.868
$ws.Range("L132").Value = 5500.0002
$ws.Range("M132").Value = -10366.868
$ws.Range("N132").Value = -10560.0002

$ws.Range("H137").Value = 6671.0713
$ws.Range("I137").Value = 4230.4116
$ws.Range("J137").Value = 10443
$ws.Range("K137").Value = 12691.2348
$ws.Range("L137").Value = 31329
$ws.Range("M137").Value = -10141.2348
$ws.Range("N137").Value = -36429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1302.2307
$ws.Range("I2").Value = 1233.75
$ws.Range("J2").Value = 1411.8
$ws.Range("K2").Value = 1233.75
$ws.Range("L2").Value = 1411.8
$ws.Range("M2").Value = -1120.75
$ws.Range("N2").Value = -1637.8

$ws.Range("H32").Value = 3475.014
$ws.Range("I32").Value = 3460.0144
$ws.Range("J32").Value = 4000
$ws.Range("K32").Value = 3460.0144
$ws.Range("L32").Value = 4000
$ws.Range("M32").Value = -3173.0144
$ws.Range("N32").Value = -4574

$ws.Range("H61").Value = 2996.1853
$ws.Range("I61").Value = 2123.45
$ws.Range("J61").Value = 5489.7144
$ws.Range("K61").Value = 2123.45
$ws.Range("L61").Value = 5489.7144
$ws.Range("M61").Value = -1911.45
$ws.Range("N61").Value = -5913.7144

$ws.Range("H110").Value = 28577.309
$ws.Range("I110").Value = 32864.09
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 32864.09
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = -30819.09
$ws.Range("N110").Value = -9090

$ws.Range("H116").Value = 1302.2307
$ws.Range("I116").Value = 1233.75
$ws.Range("J116").Value = 1411.8
$ws.Range("K116").Value = 1233.75
$ws.Range("L116").Value = 1411.8
$ws.Range("M116").Value = 1060.25
$ws.Range("N116").Value = -5999.8

$ws.Range("H122").Value = 2131.8125
$ws.Range("I122").Value = 2003.7693
$ws.Range("J122").Value = 2686.6667
$ws.Range("K122").Value = 6011.3079
$ws.Range("L122").Value = 8060.000100000001
$ws.Range("M122").Value = -3561.3079
$ws.Range("N122").Value = -12960.0001

$ws.Range("H132").Value = 20326.27
$ws.Range("I132").Value = 4902.9014
$ws.Range("J132").Value = 176487.88
$ws.Range("K132").Value = 14708.7042
$ws.Range("L132").Value = 529463.64
$ws.Range("M132").Value = -12178.7042
$ws.Range("N132").Value = -534523.64

$ws.Range("H136").Value = 2996.1853
$ws.Range("I136").Value = 2123.45
$ws.Range("J136").Value = 5489.7144
$ws.Range("K136").Value = 6370.349999999999
$ws.Range("L136").Value = 16469.1432
$ws.Range("M136").Value = -3820.349999999999
$ws.Range("N136").Value = -21569.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1302.2307
$ws.Range("I3").Value = 1233.75
$ws.Range("J3").Value = 1411.8
$ws.Range("K3").Value = 1233.75
$ws.Range("L3").Value = 1411.8
$ws.Range("M3").Value = -1119.75
$ws.Range("N3").Value = -1639.8

$ws.Range("H74").Value = 26289
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 26289
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 26289
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -28161

$ws.Range("H77").Value = 26289
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 26289
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 78867
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -88227

$ws.Range("H105").Value = 6633.15
$ws.Range("I105").Value = 6719.1055
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 6719.1055
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -4972.1055
$ws.Range("N105").Value = -8494

$ws.Range("H107").Value = 3667.6296
$ws.Range("I107").Value = 3477.5293
$ws.Range("J107").Value = 3990.8
$ws.Range("K107").Value = 3477.5293
$ws.Range("L107").Value = 3990.8
$ws.Range("M107").Value = -1557.5293
$ws.Range("N107").Value = -7830.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4935.231
$ws.Range("I31").Value = 3124.2068
$ws.Range("J31").Value = 7218.696
$ws.Range("K31").Value = 3124.2068
$ws.Range("L31").Value = 7218.696
$ws.Range("M31").Value = -2829.2068
$ws.Range("N31").Value = -7808.696

$ws.Range("H34").Value = 4935.231
$ws.Range("I34").Value = 3124.2068
$ws.Range("J34").Value = 7218.696
$ws.Range("K34").Value = 3124.2068
$ws.Range("L34").Value = 7218.696
$ws.Range("M34").Value = -2922.2068
$ws.Range("N34").Value = -7622.696

$ws.Range("H103").Value = 16841
$ws.Range("I103").Value = 16841
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 16841
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -15669

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 8529.666999999999
$ws.Range("I14").Value = 8529.666999999999
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 25589.001
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -25416.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 358.82352
$ws.Range("I2").Value = 102.2
$ws.Range("J2").Value = 465.75
$ws.Range("K2").Value = 102.2
$ws.Range("L2").Value = 465.75
$ws.Range("M2").Value = 10.8
$ws.Range("N2").Value = -691.75

$ws.Range("H113").Value = 4003.2778
$ws.Range("I113").Value = 3504.3
$ws.Range("J113").Value = 4627
$ws.Range("K113").Value = 3504.3
$ws.Range("L113").Value = 4627
$ws.Range("M113").Value = -1334.3
$ws.Range("N113").Value = -8967

$ws.Range("H122").Value = 502.5
$ws.Range("I122").Value = 520.5789
$ws.Range("J122").Value = 159
$ws.Range("K122").Value = 1561.7367
$ws.Range("L122").Value = 477
$ws.Range("M122").Value = 888.2633000000001
$ws.Range("N122").Value = -5377

$ws.Range("H132").Value = 1249.8948
$ws.Range("I132").Value = 1252.6666
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 3757.9998
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -1227.9998
$ws.Range("N132").Value = -8660

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5807.8
$ws.Range("I7").Value = 6928
$ws.Range("J7").Value = 5327.7144
$ws.Range("K7").Value = 6928
$ws.Range("L7").Value = 5327.7144
$ws.Range("M7").Value = -6816
$ws.Range("N7").Value = -5551.7144

$ws.Range("H16").Value = 3596.65
$ws.Range("I16").Value = 3473.7693
$ws.Range("J16").Value = 3824.8572
$ws.Range("K16").Value = 3473.7693
$ws.Range("L16").Value = 3824.8572
$ws.Range("M16").Value = -3303.7693
$ws.Range("N16").Value = -4164.8572

$ws.Range("H61").Value = 22729698
$ws.Range("I61").Value = 27780372
$ws.Range("J61").Value = 1670
$ws.Range("K61").Value = 27780372
$ws.Range("L61").Value = 1670
$ws.Range("M61").Value = -27780170
$ws.Range("N61").Value = -2074

$ws.Range("H113").Value = 22729698
$ws.Range("I113").Value = 27780372
$ws.Range("J113").Value = 1670
$ws.Range("K113").Value = 27780372
$ws.Range("L113").Value = 1670
$ws.Range("M113").Value = -27778202
$ws.Range("N113").Value = -6010

$ws.Range("H122").Value = 4050.818
$ws.Range("I122").Value = 3551.8572
$ws.Range("J122").Value = 4924
$ws.Range("K122").Value = 10655.5716
$ws.Range("L122").Value = 14772
$ws.Range("M122").Value = -8205.571599999999
$ws.Range("N122").Value = -19672

$ws.Range("H126").Value = 5807.8
$ws.Range("I126").Value = 6928
$ws.Range("J126").Value = 5327.7144
$ws.Range("K126").Value = 20784
$ws.Range("L126").Value = 15983.1432
$ws.Range("M126").Value = -18314
$ws.Range("N126").Value = -20923.1432

$ws.Range("H132").Value = 7078.3706
$ws.Range("I132").Value = 3558.6191
$ws.Range("J132").Value = 19397.5
$ws.Range("K132").Value = 10675.8573
$ws.Range("L132").Value = 58192.5
$ws.Range("M132").Value = -8145.8573
$ws.Range("N132").Value = -63252.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1096.5217
$ws.Range("I113").Value = 514.4
$ws.Range("J113").Value = 2188
$ws.Range("K113").Value = 1543.2
$ws.Range("L113").Value = 6564
$ws.Range("M113").Value = 626.8000000000002
$ws.Range("N113").Value = -10904
